$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a changed-date value that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 135).
for ($row = 2; $row -le 135; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
